# Commit: "fully input access db"
# The "namn" (name) column H is being removed/cleared across the whole
# used range (header in H1 plus all 35 data rows, H2:H36). Saving the
# workbook afterwards lets Excel recompact the shared-string table and
# drop the now-unused strings, matching the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1:H36").Value = ""
